$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 54 gets the "log" style (A: plain/no-color) and "description" style (C/D: wrap text),
# matching the existing pattern used by the other cwl_log_*/cwl_warn_* row pairs.
[void]$ws.Range("A50").Copy()
[void]$ws.Range("A54").PasteSpecial(-4122)   # xlPasteFormats
[void]$ws.Range("C50").Copy()
[void]$ws.Range("C54").PasteSpecial(-4122)
[void]$ws.Range("D54").PasteSpecial(-4122)

# Row 55 keeps the existing "warn" style on A (already correct) and needs the
# wrap-text description style copied onto C/D.
[void]$ws.Range("C51").Copy()
[void]$ws.Range("C55").PasteSpecial(-4122)
[void]$ws.Range("D55").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Row 54: cwl_log_custom_trait / qualified custom trait id: {0}, type: {1}
$ws.Range("A54").Value = "cwl_log_custom_trait"
$ws.Range("C54").Value = "qualified custom trait id: {0}, type: {1}"
$ws.Range("D54").Value = "qualified custom trait id: {0}, type: {1}"

# Row 55: cwl_warn_qualify_trait / failed to qualify custom trait id: {0}, on card: {1}
$ws.Range("A55").Value = "cwl_warn_qualify_trait"
$ws.Range("C55").Value = "failed to qualify custom trait id: {0}, on card: {1}"
$ws.Range("D55").Value = "failed to qualify custom trait id: {0}, on card: {1}"

# Update selection to match the new active range
[void]$ws.Range("A54:D55").Select()
